$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "56÷7=8, 0"
$t.Cell(1, 2).Range.Text  = "49÷8=6, 1"
$t.Cell(1, 3).Range.Text  = "12÷5=2, 2"
$t.Cell(1, 4).Range.Text  = "69÷4=17, 1"
$t.Cell(1, 5).Range.Text  = "15÷3=5, 0"

$t.Cell(5, 1).Range.Text  = "68÷7=9, 5"
$t.Cell(5, 2).Range.Text  = "37÷3=12, 1"
$t.Cell(5, 3).Range.Text  = "91÷5=18, 1"
$t.Cell(5, 4).Range.Text  = "12÷5=2, 2"
$t.Cell(5, 5).Range.Text  = "63÷5=12, 3"

$t.Cell(9, 1).Range.Text  = "20÷2=10, 0"
$t.Cell(9, 2).Range.Text  = "89÷6=14, 5"
$t.Cell(9, 3).Range.Text  = "92÷3=30, 2"
$t.Cell(9, 4).Range.Text  = "90÷3=30, 0"
$t.Cell(9, 5).Range.Text  = "41÷4=10, 1"

$t.Cell(13, 1).Range.Text = "25÷2=12, 1"
$t.Cell(13, 2).Range.Text = "30÷2=15, 0"
$t.Cell(13, 3).Range.Text = "37÷5=7, 2"
$t.Cell(13, 4).Range.Text = "74÷5=14, 4"
$t.Cell(13, 5).Range.Text = "43÷4=10, 3"

$t.Cell(17, 1).Range.Text = "67÷2=33, 1"
$t.Cell(17, 2).Range.Text = "88÷3=29, 1"
$t.Cell(17, 3).Range.Text = "50÷2=25, 0"
$t.Cell(17, 4).Range.Text = "93÷2=46, 1"
$t.Cell(17, 5).Range.Text = "35÷3=11, 2"
